$d = $word.ActiveDocument

# Locate the paragraph that holds the M2Doc field
# ( {begin} m:'doc.html'.fromHTMLURI() {end} ) so we can rewrite it as
# plain-text runs "{", "m", ":", "'", "doc.html", "'.fromHTMLURI()", "}".
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Fields.Count -gt 0) {
        $target = $candidate
    }
}

if ($target -ne $null) {
    # Recover the paragraph's original opening-tag attributes (paraId,
    # rsid*, ...) so the rewrite only touches the runs, as in the diff.
    $pTag = "<w:p>"
    $fullXml = $d.WordOpenXML
    $anchor = $fullXml.IndexOf('fldChar w:fldCharType="begin"')
    if ($anchor -ge 0) {
        $head = $fullXml.Substring(0, $anchor)
        $pStart = $head.LastIndexOf("<w:p ")
        if ($pStart -ge 0) {
            $pEnd = $fullXml.IndexOf(">", $pStart)
            if ($pEnd -gt $pStart) {
                $pTag = $fullXml.Substring($pStart, $pEnd - $pStart + 1)
            }
        }
    }
    if ($pTag.EndsWith("/>")) {
        $pTag = $pTag.Substring(0, $pTag.Length - 2) + ">"
    }

    $newParaXml = $pTag +
        "<w:r><w:t>{</w:t></w:r>" +
        "<w:r><w:t>m</w:t></w:r>" +
        "<w:r><w:t>:</w:t></w:r>" +
        "<w:r><w:t>'</w:t></w:r>" +
        "<w:r><w:t>doc.html</w:t></w:r>" +
        "<w:r><w:t>'.fromHTMLURI()</w:t></w:r>" +
        "<w:r><w:t xml:space='preserve'>}</w:t></w:r>" +
        "</w:p>"

    $target.Range.InsertXML($newParaXml)
}
